# The workbook tracks LDLC smartphone prices over time.
# Each scrape run appends a new timestamped price column just before the
# two trailing "nom" (product name) and "url_produit" columns.
#
# Before this edit the trailing columns are:
#   ... | FL = last price snapshot | FM = nom | FN = url_produit
# After this edit a brand-new snapshot column is inserted before "nom",
# carrying forward the most recent known price (or staying blank for rows
# whose price tracking had already stopped), and "nom"/"url_produit" shift
# one column to the right:
#   ... | FL = last price snapshot | FM = new snapshot | FN = nom | FO = url_produit

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new column right before the "nom" column (FM).
# This shifts the existing "nom" -> FN and "url_produit" -> FO, and Excel
# automatically extends the sheet dimension/used range to include it.
$ws.Range("FM1").EntireColumn.Insert()

# Header for the freshly inserted column: the new scrape timestamp.
$ws.Range("FM1").Value = "2026-02-04 18:33:31"

# Populate the new column's data rows: carry forward the previous price
# snapshot (column FL) wherever it had a value; leave it blank for rows
# where price tracking had already stopped (FL empty).
$lastRow = 208
$priceCol = 168   # column FL
$newCol = 169     # column FM (newly inserted)

for ($r = 2; $r -le $lastRow; $r++) {
  $flVal = $ws.Cells.Item($r, $priceCol).Value()
  if ($flVal -ne $null -and $flVal -ne "") {
    $ws.Cells.Item($r, $newCol).Value = $flVal
  }
}
